$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Release Date (RD)")

$ws.Range("A2").Value = 253393
$ws.Range("B2").Value = 45912.58333333334
$ws.Range("C2").Value = 0

$ws.Range("A3").Value = 253392
$ws.Range("B3").Value = 45911.58333333334
$ws.Range("B3").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("C3").Value = 0
